$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("CompleteSVM")

$sheetsToAdd = @(
    @{ Name = "CompleteSVM1"; F1 = 0.0459498728659644; Accuracy = 0.1012 },
    @{ Name = "CompleteSVM2"; F1 = 0.4821801076932929; Accuracy = 0.7312 },
    @{ Name = "CompleteSVM3"; F1 = 0.4821801076932929; Accuracy = 0.7312 },
    @{ Name = "CompleteSVM4"; F1 = 0.9831150166659732; Accuracy = 0.9844000000000001 }
)

foreach ($info in $sheetsToAdd) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $info.Name

    $template.Range("B1:C1").Copy()
    $ws.Range("B1:C1").PasteSpecial(-4122)

    $template.Range("A2").Copy()
    $ws.Range("A2").PasteSpecial(-4122)

    $ws.Range("B1").Value = "f1_score"
    $ws.Range("C1").Value = "accuracy"

    $ws.Range("A2").Value = 0
    $ws.Range("B2").Value = $info.F1
    $ws.Range("C2").Value = $info.Accuracy
}
